# Update column G ("K" - strikeouts) values to the regenerated s_vals.
# Commit: "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 7
    3  = 8
    4  = 8
    5  = 9
    6  = 4
    7  = 8
    8  = 8
    9  = 10
    10 = 14
    11 = 4
    12 = 11
    13 = 5
    14 = 6
    15 = 8
    16 = 8
    17 = 9
    18 = 6
    19 = 12
    20 = 11
    21 = 7
    22 = 6
    23 = 8
    24 = 12
    25 = 7
    26 = 3
    27 = 9
    28 = 10
    29 = 11
    30 = 5
    31 = 9
    32 = 10
    33 = 6
    34 = 9
    35 = 3
    36 = 6
    37 = 3
    38 = 5
    39 = 5
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
